$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9835058450698853
$ws.Range("B1").Value = 1.978021025657654
$ws.Range("C1").Value = 8.690386772155762
$ws.Range("D1").Value = 2.802839517593384
$ws.Range("E1").Value = 1.426483392715454
